# Update the "GLDW_DataSources" worksheet (Sheet1):
#  - Row 4 (SimulatedBuoys) description text is extended with a note about the NDBC data stream.
#  - Row 5 changes from the retired "TDSManual / TDSData" (Total Dissolved Solids) data source
#    to the new "SalinityProject / SalinityData" data source, including an updated description
#    and an updated link pointing at the Salinity entry form.
#  - The active selection moves to F5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = "This is simulated data demonstrating how data from a series of buoys on Lake Erie could be graphed, analyzed and trigger alerts. The NDBC data stream features actual buoy data."

$ws.Range("B5").Value = "SalinityData"
$ws.Range("A5").Value = "SalinityProject"
$ws.Range("E5").Value = "This demonstrates how specific conductivity data manually entered from a smartphone can be validated prior to incorporating the data in the GLDW  repository."
$ws.Range("F5").Value = "<a href='http://sources.gldw.org:40402/vdab'>Container</a><hr><a href='http://sources.gldw.org:40402/vdab/views/enterSalinity'>EnterData</a>"

# Row 4 grows an extra wrapped line, so its autofit height increases from 60 to 75.
$ws.Rows("4").RowHeight = 75

# Move/select the active cell to F5 (matches the saved selection in the sheet view).
$ws.Range("F5").Select()
